$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: Alps SKRPABE010 key switch count 68 -> 63
$ws.Range("F19").Value = "63"
$ws.Range("G19").Value = "63キー / 4.2×3.2mm SMD / JLCPCB在庫54355個確認済み"

# Row 20: 1N4148W diode count 66 -> 63
$ws.Range("F20").Value = "63"
$ws.Range("G20").Value = "キーマトリクス用（63キー分）"
